$d = $word.ActiveDocument

$replacements = @(
    @{old = "697×2="; new = "565×6="},
    @{old = "155×4="; new = "565×2="},
    @{old = "424×4="; new = "670×8="},
    @{old = "213×5="; new = "839×5="},
    @{old = "836×4="; new = "563×6="},
    @{old = "700×5="; new = "849×4="},
    @{old = "294×5="; new = "543×4="},
    @{old = "498×2="; new = "194×2="},
    @{old = "270×7="; new = "513×5="},
    @{old = "355×5="; new = "730×4="},
    @{old = "731×3="; new = "728×4="},
    @{old = "340×9="; new = "537×3="},
    @{old = "920×6="; new = "803×8="},
    @{old = "113×8="; new = "694×9="},
    @{old = "803×7="; new = "464×2="},
    @{old = "314×7="; new = "572×5="},
    @{old = "638×6="; new = "759×3="},
    @{old = "778×9="; new = "653×3="},
    @{old = "386×4="; new = "554×9="},
    @{old = "318×5="; new = "534×6="},
    @{old = "726×9="; new = "997×7="},
    @{old = "468×4="; new = "923×7="},
    @{old = "887×6="; new = "562×3="},
    @{old = "500×4="; new = "248×8="},
    @{old = "481×5="; new = "396×4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
